$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D=44495; H="Sin especificar"; I="Primera"; J=200; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=3; D=44194; H="Sin especificar"; I="Extra"; J=120; K=3500; L=3500; M=3500; N="$/unidad"; O="Región de O'Higgins"; P=3500 },
    @{ Row=4; D=44194; H="Sin especificar"; I="Primera"; J=200; K=3000; L=3000; M=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000 },
    @{ Row=5; D=44312; H="Sin especificar"; I="Primera"; J=180; K=2500; L=2500; M=2500; N="$/unidad"; O="Perú"; P=2500 },
    @{ Row=6; D=44167; H="Sin especificar"; I="Primera"; J=400; K=5000; L=5000; M=5000; N="$/unidad"; O="Región de O'Higgins"; P=5000 },
    @{ Row=7; D=44167; H="Sin especificar"; I="Segunda"; J=560; K=3000; L=3000; M=3000; N="$/unidad"; O="Región de O'Higgins"; P=3000 },
    @{ Row=8; D=44167; H="Sin especificar"; I="Tercera"; J=450; K=2000; L=2000; M=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000 },
    @{ Row=9; D=44477; H="Sin especificar"; I="Primera"; J=80; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=10; D=44483; H="Sin especificar"; I="Primera"; J=120; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=11; D=44217; H="Sin especificar"; I="Extra"; J=400; K=2500; L=2500; M=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500 },
    @{ Row=12; D=44217; H="Sin especificar"; I="Primera"; J=280; K=2000; L=2000; M=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000 },
    @{ Row=13; D=44504; H="Sin especificar"; I="Primera"; J=200; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=14; D=44223; H="Americana O Klondike"; I="Extra"; J=340; K=2500; L=2500; M=2500; N="$/unidad"; O="Región de O'Higgins"; P=2500 },
    @{ Row=15; D=44223; H="Americana O Klondike"; I="Primera"; J=400; K=2000; L=2000; M=2000; N="$/unidad"; O="Región de O'Higgins"; P=2000 },
    @{ Row=16; D=44223; H="Americana O Klondike"; I="Segunda"; J=300; K=1500; L=1500; M=1500; N="$/unidad"; O="Región de O'Higgins"; P=1500 },
    @{ Row=17; D=44223; H="Americana O Klondike"; I="Tercera"; J=160; K=1000; L=1000; M=1000; N="$/unidad"; O="Región de O'Higgins"; P=1000 },
    @{ Row=18; D=44488; H="Sin especificar"; I="Primera"; J=150; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=19; D=44510; H="Sin especificar"; I="Primera"; J=250; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=20; D=44497; H="Sin especificar"; I="Primera"; J=250; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=21; D=44491; H="Sin especificar"; I="Primera"; J=150; K=800; L=800; M=800; N="$/kilo (volumen en unidades)"; O="Perú"; P=800 },
    @{ Row=22; D=44305; H="Sin especificar"; I="Primera"; J=100; K=2500; L=2500; M=2500; N="$/unidad"; O="Perú"; P=2500 },
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
}
